$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "51.684.00"
Set-TextValue $ws.Range("D3") "2.783.50"
Set-TextValue $ws.Range("E3") "  -0.25%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.03%  "
Set-TextValue $ws.Range("D5") "352.07"
Set-TextValue $ws.Range("E5") "  -2.02%  "
Set-TextValue $ws.Range("D6") "108.83"
Set-TextValue $ws.Range("E6") "  -0.56%  "
Set-TextValue $ws.Range("D7") "0.547"
Set-TextValue $ws.Range("E7") "  -3.02%  "
Set-TextValue $ws.Range("E8") "  +0.10%  "
Set-TextValue $ws.Range("E9") "  +1.91%  "
Set-TextValue $ws.Range("E10") "  -0.49%  "
Set-TextValue $ws.Range("E11") "  +2.56%  "
Set-TextValue $ws.Range("D12") "20.13"
Set-TextValue $ws.Range("E12") "  +3.50%  "
Set-TextValue $ws.Range("E13") "  -2.05%  "
Set-TextValue $ws.Range("D14") "7.68"
Set-TextValue $ws.Range("E14") "  +1.30%  "
Set-TextValue $ws.Range("D15") "3.221.34"
Set-TextValue $ws.Range("E15") "  -0.14%  "
Set-TextValue $ws.Range("D16") "2.794.03"
Set-TextValue $ws.Range("E16") "  +0.00%  "
Set-TextValue $ws.Range("D17") "0.929"
Set-TextValue $ws.Range("E17") "  -2.12%  "
Set-TextValue $ws.Range("D18") "51.577.05"
Set-TextValue $ws.Range("E18") "  -0.58%  "
Set-TextValue $ws.Range("D19") "7.73"
Set-TextValue $ws.Range("E19") "  +4.50%  "
Set-TextValue $ws.Range("D20") "3.14"
Set-TextValue $ws.Range("E20") "  +0.11%  "
Set-TextValue $ws.Range("D21") "13.17"
Set-TextValue $ws.Range("E21") "  +1.33%  "
Set-TextValue $ws.Range("E22") "  -1.64%  "
Set-TextValue $ws.Range("D23") "69.90"
Set-TextValue $ws.Range("E23") "  -0.51%  "
Set-TextValue $ws.Range("D24") "267.20"
Set-TextValue $ws.Range("E24") "  -2.65%  "
Set-TextValue $ws.Range("D25") "2.74"
Set-TextValue $ws.Range("E25") "  -0.36%  "
Set-TextValue $ws.Range("E26") "  -2.23%  "
Set-TextValue $ws.Range("E27") "  -0.07%  "
Set-TextValue $ws.Range("D28") "0.164"
Set-TextValue $ws.Range("E28") "  +13.00%  "
Set-TextValue $ws.Range("D29") "10.23"
Set-TextValue $ws.Range("E29") "  +0.43%  "
Set-TextValue $ws.Range("D30") "37.17"
Set-TextValue $ws.Range("E30") "  +7.39%  "
Set-TextValue $ws.Range("D31") "2.23"
Set-TextValue $ws.Range("E31") "  -2.12%  "
Set-TextValue $ws.Range("E32") "  +8.18%  "
Set-TextValue $ws.Range("E33") "  +0.11%  "
Set-TextValue $ws.Range("D34") "0.0452"
Set-TextValue $ws.Range("E34") "  -2.00%  "
Set-TextValue $ws.Range("D35") "5.63"
Set-TextValue $ws.Range("E35") "  +6.18%  "
Set-TextValue $ws.Range("D36") "0.0830"
Set-TextValue $ws.Range("E36") "  -1.92%  "
Set-TextValue $ws.Range("D37") "0.999"
Set-TextValue $ws.Range("E37") "  -0.07%  "
Set-TextValue $ws.Range("E38") "  +2.27%  "
Set-TextValue $ws.Range("E39") "  -2.36%  "
Set-TextValue $ws.Range("E40") "  -1.54%  "
Set-TextValue $ws.Range("E41") "  -1.02%  "
Set-TextValue $ws.Range("E42") "  -0.70%  "
Set-TextValue $ws.Range("E43") "  -1.66%  "
Set-TextValue $ws.Range("D44") "22.10"
Set-TextValue $ws.Range("E44") "  +0.25%  "
Set-TextValue $ws.Range("E45") "  -2.76%  "
Set-TextValue $ws.Range("D46") "2.124.70"
Set-TextValue $ws.Range("E46") "  +2.32%  "
Set-TextValue $ws.Range("E47") "  +1.43%  "
Set-TextValue $ws.Range("E48") "  +6.11%  "
Set-TextValue $ws.Range("E49") "  -5.14%  "
Set-TextValue $ws.Range("D50") "0.907"
Set-TextValue $ws.Range("E50") "  -3.00%  "
Set-TextValue $ws.Range("D51") "1.34"
Set-TextValue $ws.Range("E51") "  +9.35%  "
